# Swap the "Nashul Aptel" / "Anshul Patel" rows (rows 3 and 4) on Sheet1,
# per "WA Beta" update, and move the active selection to A3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture current values for row 3 and row 4 (columns A and B)
$a3 = $ws.Range("A3").Value2
$b3 = $ws.Range("B3").Value2
$a4 = $ws.Range("A4").Value2
$b4 = $ws.Range("B4").Value2

# Swap the rows' contents
$ws.Range("A3").Value = $a4
$ws.Range("B3").Value = $b4
$ws.Range("A4").Value = $a3
$ws.Range("B4").Value = $b3

# Update the active selection to A3 (was A10)
$ws.Range("A3").Select()
